$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 (shifts TORSAMOLEX...and everything below down by one,
# including the totals row and the footer row).
$ws.Rows("22:22").Insert()

# Copy the formatting (styles, borders, number formats) of the row that used to be at 22
# (now at 23) into the freshly inserted blank row so it matches the rest of the table.
$ws.Range("A23:N23").Copy()
$ws.Range("A22:N22").PasteSpecial(-4122)
$ws.Rows("22:22").RowHeight = $ws.Rows("23:23").RowHeight

# Re-create the merged cells for the new row (B:G, H:K, L:M), matching every other
# item row in the table.
$ws.Range("B22:G22").Merge()
$ws.Range("H22:K22").Merge()
$ws.Range("L22:M22").Merge()

# Fill in the new item: "RELAXON 30 CAP", inserted alphabetically between
# "PANADOL ADVANCE 500 MG 48 TABLETS" (row 21) and "TORSAMOLEX 20MG 20 TABS" (now row 23).
$ws.Range("A22").Value = 19
$ws.Range("B22").Value = "RELAXON 30 CAP"
$ws.Range("H22").Value = "1:2"
$ws.Range("L22").Value = 23
$ws.Range("N22").Value = "0:0"

# Renumber the running index column (A) for every row that shifted down because of the insert.
For ($r = 23; $r -le 34; $r++) {
    $ws.Range("A$r").Value = $r - 3
}

# Update the grand-total cell (now on row 35) to reflect the new row's price value.
$ws.Range("K35").Value = 1012.32
